# TC_52_Verify_Devices_Order_in_Dropped_Gallery.xlsx
# "Updated test data for power calculation text and other base change related
#  test cases"
#
# The "Device Order Gallery Expansion" sheet had a stale "801 CHEx IS" test
# row (Device Index 2) in its Detectors block; that row is removed and the
# remaining rows shift up, so what was row 13 (Device Index 3 / "801 F")
# becomes row 12. The previously-active sheet/selection ("Device Order
# Context Menu"!B11) reverts to the first sheet, now selected at C11.

$wb = $excel.ActiveWorkbook

$wsGallery = $wb.Worksheets.Item("Device Order Gallery Expansion")
$wsContext = $wb.Worksheets.Item("Device Order Context Menu")

# Remove the obsolete "801 CHEx IS" row (row 12); rows below shift up so the
# former row 13 ("801 F", index 3) becomes the new row 12.
$wsGallery.Rows.Item(12).Delete()

# The gallery sheet becomes the active sheet/tab again, with C11 selected.
$wsGallery.Activate()
$wsGallery.Range("C11").Select()

# Context menu sheet keeps its own prior selection (B11), it just stops
# being the active tab once the gallery sheet is activated above.
$wsContext.Range("B11").Select()
$wsGallery.Activate()
